{"js": "// The document's first paragraph carries the AFFARS topic-id marker.\nconst body = context.document.body;\nconst firstPara = body.paragraphs.getFirst();\nfirstPara.load(\"text\");\nawait context.sync();\n\n// --- 1) Paragraph border (space-only, 5pt on every side) ---------------\n// Office.js's public Paragraph API has no \"distance from text\" property\n// for paragraph borders (Word.Border only exposes color/lineStyle/width),\n// so we reach through the low-level OM bridge (the same dispatch the\n// Word object model / VBA \"ParagraphFormat.Borders.DistanceFromXxx\"\n// properties use) to set it without perturbing anything else in the run.\nconst borders = firstPara.borders;\nborders._omSet(\"DistanceFromTop\", 5, \"Borders\");\nborders._omSet(\"DistanceFromLeft\", 5, \"Borders\");\nborders._omSet(\"DistanceFromBottom\", 5, \"Borders\");\nborders._omSet(\"DistanceFromRight\", 5, \"Borders\");\n\n// --- 2) Left indent: 120 twips (6pt) -> 225 twips (11.25pt) -------------\nfirstPara.leftIndent = 11.25;\n\n// --- 3) Marker text update + collapse the stray trailing-space run -----\n// The paragraph currently holds two runs (\"**ID__...__ID**\" and a lone\n// trailing \" \"); replacing the whole paragraph's text with the new\n// marker collapses it back down to a single run with no trailing space.\nfirstPara.insertText(\"**ID__AFFARS_SMC_PGI_5322__ID**\", Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# First paragraph of the body holds the AFFARS topic id marker.\n$p = $d.Paragraphs(1)\n\n# Add a paragraph border (all four sides, space-only / no visible line)\n# with 5pt spacing on each side, and bump the left indent from 6pt (120\n# twips) to 11.25pt (225 twips).\n$pFmt = $p.Range.ParagraphFormat\n$pFmt.Borders.DistanceFromTop = 5\n$pFmt.Borders.DistanceFromLeft = 5\n$pFmt.Borders.DistanceFromBottom = 5\n$pFmt.Borders.DistanceFromRight = 5\n$pFmt.LeftIndent = 11.25\n\n# Replace the paragraph's text (excluding the trailing paragraph mark) with\n# the updated marker, collapsing the old two runs (marker + trailing\n# space) down into a single run with no trailing space.\n$r = $p.Range\n$textRange = $d.Range($r.Start, $r.End - 1)\n$textRange.Text = \"**ID__AFFARS_SMC_PGI_5322__ID**\"\n"}
